$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 2024
$ws.Range("B6").Value = "DEC"
$ws.Range("C6").Value = "31/12-01/12"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "Akurana"
$ws.Range("F6").Value = "✓"
$ws.Range("G6").Value = "✓"
$ws.Range("H6").Value = "✓"
$ws.Range("I6").Value = "✓"
$ws.Range("J6").Value = "✓"
$ws.Range("K6").Value = "✓"
$ws.Range("L6").Value = "✓"
$ws.Range("M6").Value = "-"
$ws.Range("N6").Value = "✓"
$ws.Range("O6").Value = "✓"
$ws.Range("P6").Value = "✓"
$ws.Range("Q6").Value = "✓"
$ws.Range("R6").Value = "✓"
$ws.Range("S6").Value = "✓"
$ws.Range("T6").Value = "✓"
$ws.Range("U6").Value = "✓"
$ws.Range("V6").Value = "✓"
$ws.Range("W6").Value = "✓"
$ws.Range("X6").Value = "✓"
$ws.Range("Y6").Value = "-"
$ws.Range("Z6").Value = "-"
$ws.Range("AA6").Value = "✓"
$ws.Range("AB6").Value = "✓"
$ws.Range("AC6").Value = "✓"
$ws.Range("AD6").Value = "✓"
$ws.Range("AE6").Value = "✓"
$ws.Range("AF6").Value = "✓"
$ws.Range("AG6").Value = "✓"
$ws.Range("AH6").Value = "-"
$ws.Range("AI6").Value = "-"
$ws.Range("AJ6").Value = "-"
